$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "2023" column (S) data ---
$ws.Range("S3").Value = 2023
$ws.Range("S4").Value = 3351.5
$ws.Range("S5").Value = 388.4
$ws.Range("S6").Value = 10593
$ws.Range("S7").Value = 7732
$ws.Range("S8").Value = 942.3
$ws.Range("S9").Value = 104.1
$ws.Range("S10").Value = 88.6
$ws.Range("S11").Value = 284.3
$ws.Range("S12").Value = 11.4
$ws.Range("S13").Value = 3094.1
$ws.Range("S14").Value = 2999.7

# Copy styling from the previous "2022" column (R) onto the new "2023" column (S)
$ws.Range("R1:R15").Copy()
$ws.Range("S1:S15").PasteSpecial(-4122)

# --- Shift the empty "blank separator" column to U ---
$ws.Range("T2:T15").Copy()
$ws.Range("U2:U15").PasteSpecial(-4122)

# --- Update the merged title cell to span through the new column ---
$ws.Range("A1:R1").UnMerge()
$ws.Range("A1:S1").Merge()

# --- Update the active selection to match the new layout ---
$ws.Range("S3:S14").Select() | Out-Null
